# Applies the "Share Product 15 Test Cases" edit:
#  - Input sheet: Groupname/AddClientMember/verify1 rows updated to new
#    "create group + attach client" test-case data (Group4108, Jhon Deer,
#    GroupAddClient/click)
#  - Output sheet: verifyGroup row now points at Group4108 and a new
#    VerifyClientCreated/Jhon Deer row is appended
#  - Selections are updated to match the saved workbook state
#
# NOTE: the cell-write order below is deliberate (not just cosmetic) -
# it reproduces the order in which new text values were first entered
# so that the shared-string table comes out in the same order as the
# original author's save.

$wb = $excel.ActiveWorkbook

$wsInput = $wb.Worksheets.Item("Input")
$wsOutput = $wb.Worksheets.Item("Output")

# 1) Input!A5:B5 - AddClientMember / Jhon Deer
$wsInput.Range("A5").Value = "AddClientMember"
$wsInput.Range("B5").Value = "Jhon Deer"

# 2) Input!A6:B6 - GroupAddClient / click
$wsInput.Range("A6").Value = "GroupAddClient"
$wsInput.Range("B6").Value = "click"

# 3) Output!A2 - VerifyClientCreated (new row, text added first)
$wsOutput.Range("A2").Value = "VerifyClientCreated"

# 4) Input!A2:B2 - Groupname / Group4108
$wsInput.Range("A2").Value = "Groupname"
$wsInput.Range("B2").Value = "Group4108"

# 5) Output!A1:B1 - verifyGroup / Group4108
$wsOutput.Range("A1").Value = "verifyGroup"
$wsOutput.Range("B1").Value = "Group4108"

# Restyle Output!B1 to match the other "value" cells on the Input sheet
# (Calibri font on the shaded fill) by copying the format from Input!B1.
$wsInput.Range("B1").Copy() | Out-Null
$wsOutput.Range("B1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$wsOutput.Range("B1").Value = "Group4108"

# 6) Output!A2:B2 - finish the new row with B2 (Jhon Deer) and apply the
# plain Arial 10, left/top aligned look used elsewhere in the workbook.
$wsOutput.Range("B2").Value = "Jhon Deer"
$rngNewRow = $wsOutput.Range("A2:B2")
$rngNewRow.HorizontalAlignment = -4131
$rngNewRow.VerticalAlignment = -4160
$rngNewRow.Font.Name = "Arial"
$rngNewRow.Font.Size = 10

# Update the active selection on the Input sheet to match the saved file
$wsInput.Activate()
$wsInput.Range("D16").Select()

# Update the active selection on the Output sheet to match the saved file
$wsOutput.Activate()
$wsOutput.Range("C11").Select()

$wsInput.Activate()
